# Update row 9 (Ano 2025) values in the faturamento_anual sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 4057571.61
$ws.Range("C9").Value = 639401.24
$ws.Range("D9").Value = 4696972.85
$ws.Range("E9").Value = 13.61304951975611
$ws.Range("F9").Value = 86.3869504802439
$ws.Range("G9").Value = -38.20505835556403
$ws.Range("H9").Value = -26.72586525961815
$ws.Range("I9").Value = 40821
$ws.Range("J9").Value = 1752
$ws.Range("K9").Value = 42573
$ws.Range("L9").Value = 29482
$ws.Range("M9").Value = 159.3166287904484
$ws.Range("N9").Value = 8.768707393224794
